$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 trailing duplicate rows (18, 19, 20)
$ws.Rows.Item(18).EntireRow.Delete()
$ws.Rows.Item(18).EntireRow.Delete()
$ws.Rows.Item(18).EntireRow.Delete()

# Rewrite the data rows (2-17) to match the trimmed/corrected dataset
# Row 2
$ws.Range("A2").Value = 405
$ws.Range("B2").Value = 'Lappeenrannan kaupunki'
$ws.Range("C2").Value = 44105
$ws.Range("D2").Value = 2200
$ws.Range("E2").Value = 67179
$ws.Range("F2").Value = '2790230-3'
$ws.Range("G2").Value = 'SINGA OY'
$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J2").Value = 4600
$ws.Range("K2").Value = 'Muu materiaali'
$ws.Range("L2").Value = 249
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 'Lapset ja nuoret'
$ws.Range("O2").Value = 222
$ws.Range("P2").Value = 'Nuorisopalvelut'
$ws.Range("Q2").Value = 2230

# Row 3
$ws.Range("A3").Value = 405
$ws.Range("B3").Value = 'Lappeenrannan kaupunki'
$ws.Range("C3").Value = 44105
$ws.Range("D3").Value = 2673
$ws.Range("E3").Value = 20908
$ws.Range("F3").Value = '1022684-7'
$ws.Range("G3").Value = 'SANEERAUS J.LIIKKA OY'
$ws.Range("H3").Value = 43
$ws.Range("I3").Value = 'Palvelujen ostot'
$ws.Range("J3").Value = 4390
$ws.Range("K3").Value = 'Rak.ja al.rak.-ja kunn.pitopal'
$ws.Range("L3").Value = 1144
$ws.Range("M3").Value = 30
$ws.Range("N3").Value = 'Kaupunkikehitys'
$ws.Range("O3").Value = 340
$ws.Range("P3").Value = 'Kadut ja ympäristö'
$ws.Range("Q3").Value = 3421

# Row 4
$ws.Range("A4").Value = 405
$ws.Range("B4").Value = 'Lappeenrannan kaupunki'
$ws.Range("C4").Value = 44115
$ws.Range("D4").Value = 3460
$ws.Range("E4").Value = 76785
$ws.Range("F4").Value = '2867588-6'
$ws.Range("G4").Value = 'SAIMAARIUM OY'
$ws.Range("H4").Value = 43
$ws.Range("I4").Value = 'Palvelujen ostot'
$ws.Range("J4").Value = 4440
$ws.Range("K4").Value = 'Koulutus- ja kulttuuripalvelut'
$ws.Range("L4").Value = 50.4
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 'Lapset ja nuoret'
$ws.Range("O4").Value = 221
$ws.Range("P4").Value = 'Lukiokoulutus'
$ws.Range("Q4").Value = 2211

# Row 5
$ws.Range("A5").Value = 405
$ws.Range("B5").Value = 'Lappeenrannan kaupunki'
$ws.Range("C5").Value = 44123
$ws.Range("D5").Value = 4185
$ws.Range("E5").Value = 34912
$ws.Range("F5").Value = '0165069-7'
$ws.Range("G5").Value = 'SAVONLINJA OY'
$ws.Range("H5").Value = 43
$ws.Range("I5").Value = 'Palvelujen ostot'
$ws.Range("J5").Value = 4420
$ws.Range("K5").Value = 'Matkustus- ja kuljetuspalvelut'
$ws.Range("L5").Value = 54
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 'Lapset ja nuoret'
$ws.Range("O5").Value = 217
$ws.Range("P5").Value = 'Perusopetus'
$ws.Range("Q5").Value = 2174

# Row 6
$ws.Range("A6").Value = 405
$ws.Range("B6").Value = 'Lappeenrannan kaupunki'
$ws.Range("C6").Value = 44121
$ws.Range("D6").Value = 4209
$ws.Range("E6").Value = 1670
$ws.Range("F6").Value = '1988068-5'
$ws.Range("G6").Value = 'SARCO OY/POHJOLA PANKKI OYJ'
$ws.Range("H6").Value = 45
$ws.Range("I6").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J6").Value = 4580
$ws.Range("K6").Value = 'Kalusto'
$ws.Range("L6").Value = 1169.82
$ws.Range("M6").Value = 41
$ws.Range("N6").Value = 'Maakunnallinen palvelutoiminta'
$ws.Range("O6").Value = 420
$ws.Range("P6").Value = 'E-K Pelastuslaitos'
$ws.Range("Q6").Value = 4231

# Row 7
$ws.Range("A7").Value = 405
$ws.Range("B7").Value = 'Lappeenrannan kaupunki'
$ws.Range("C7").Value = 44133
$ws.Range("D7").Value = 5474
$ws.Range("E7").Value = 34912
$ws.Range("F7").Value = '0165069-7'
$ws.Range("G7").Value = 'SAVONLINJA OY'
$ws.Range("H7").Value = 43
$ws.Range("I7").Value = 'Palvelujen ostot'
$ws.Range("J7").Value = 4420
$ws.Range("K7").Value = 'Matkustus- ja kuljetuspalvelut'
$ws.Range("L7").Value = 24
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 'Kaupunkikehitys'
$ws.Range("O7").Value = 340
$ws.Range("P7").Value = 'Kadut ja ympäristö'
$ws.Range("Q7").Value = 3400

# Row 8
$ws.Range("A8").Value = 405
$ws.Range("B8").Value = 'Lappeenrannan kaupunki'
$ws.Range("C8").Value = 44133
$ws.Range("D8").Value = 5566
$ws.Range("E8").Value = 34912
$ws.Range("F8").Value = '0165069-7'
$ws.Range("G8").Value = 'SAVONLINJA OY'
$ws.Range("H8").Value = 43
$ws.Range("I8").Value = 'Palvelujen ostot'
$ws.Range("J8").Value = 4420
$ws.Range("K8").Value = 'Matkustus- ja kuljetuspalvelut'
$ws.Range("L8").Value = 31142.04
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 'Kaupunkikehitys'
$ws.Range("O8").Value = 340
$ws.Range("P8").Value = 'Kadut ja ympäristö'
$ws.Range("Q8").Value = 3401

# Row 9
$ws.Range("A9").Value = 405
$ws.Range("B9").Value = 'Lappeenrannan kaupunki'
$ws.Range("C9").Value = 44133
$ws.Range("D9").Value = 5566
$ws.Range("E9").Value = 34912
$ws.Range("F9").Value = '0165069-7'
$ws.Range("G9").Value = 'SAVONLINJA OY'
$ws.Range("H9").Value = 43
$ws.Range("I9").Value = 'Palvelujen ostot'
$ws.Range("J9").Value = 4420
$ws.Range("K9").Value = 'Matkustus- ja kuljetuspalvelut'
$ws.Range("L9").Value = 33
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 'Kaupunkikehitys'
$ws.Range("O9").Value = 340
$ws.Range("P9").Value = 'Kadut ja ympäristö'
$ws.Range("Q9").Value = 3400

# Row 10
$ws.Range("A10").Value = 405
$ws.Range("B10").Value = 'Lappeenrannan kaupunki'
$ws.Range("C10").Value = 44121
$ws.Range("D10").Value = 4209
$ws.Range("E10").Value = 1670
$ws.Range("F10").Value = '1988068-5'
$ws.Range("G10").Value = 'SARCO OY/POHJOLA PANKKI OYJ'
$ws.Range("H10").Value = 45
$ws.Range("I10").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J10").Value = 4580
$ws.Range("K10").Value = 'Kalusto'
$ws.Range("L10").Value = 21
$ws.Range("M10").Value = 41
$ws.Range("N10").Value = 'Maakunnallinen palvelutoiminta'
$ws.Range("O10").Value = 420
$ws.Range("P10").Value = 'E-K Pelastuslaitos'
$ws.Range("Q10").Value = 4231

# Row 11
$ws.Range("A11").Value = 405
$ws.Range("B11").Value = 'Lappeenrannan kaupunki'
$ws.Range("C11").Value = 44121
$ws.Range("D11").Value = 4209
$ws.Range("E11").Value = 1670
$ws.Range("F11").Value = '1988068-5'
$ws.Range("G11").Value = 'SARCO OY/POHJOLA PANKKI OYJ'
$ws.Range("H11").Value = 45
$ws.Range("I11").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J11").Value = 4580
$ws.Range("K11").Value = 'Kalusto'
$ws.Range("L11").Value = 222
$ws.Range("M11").Value = 41
$ws.Range("N11").Value = 'Maakunnallinen palvelutoiminta'
$ws.Range("O11").Value = 420
$ws.Range("P11").Value = 'E-K Pelastuslaitos'
$ws.Range("Q11").Value = 4231

# Row 12
$ws.Range("A12").Value = 405
$ws.Range("B12").Value = 'Lappeenrannan kaupunki'
$ws.Range("C12").Value = 44105
$ws.Range("D12").Value = 2673
$ws.Range("E12").Value = 20908
$ws.Range("F12").Value = '1022684-7'
$ws.Range("G12").Value = 'SANEERAUS J.LIIKKA OY'
$ws.Range("H12").Value = 43
$ws.Range("I12").Value = 'Palvelujen ostot'
$ws.Range("J12").Value = 4390
$ws.Range("K12").Value = 'Rak.ja al.rak.-ja kunn.pitopal'
$ws.Range("L12").Value = 1144
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 'Kaupunkikehitys'
$ws.Range("O12").Value = 340
$ws.Range("P12").Value = 'Kadut ja ympäristö'
$ws.Range("Q12").Value = 3421

# Row 13
$ws.Range("A13").Value = 405
$ws.Range("B13").Value = 'Lappeenrannan kaupunki'
$ws.Range("C13").Value = 44105
$ws.Range("D13").Value = 2673
$ws.Range("E13").Value = 20908
$ws.Range("F13").Value = '1022684-7'
$ws.Range("G13").Value = 'SANEERAUS J.LIIKKA OY'
$ws.Range("H13").Value = 43
$ws.Range("I13").Value = 'Palvelujen ostot'
$ws.Range("J13").Value = 4390
$ws.Range("K13").Value = 'Rak.ja al.rak.-ja kunn.pitopal'
$ws.Range("L13").Value = 222
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 'Kaupunkikehitys'
$ws.Range("O13").Value = 340
$ws.Range("P13").Value = 'Kadut ja ympäristö'
$ws.Range("Q13").Value = 3421

# Row 14
$ws.Range("A14").Value = 405
$ws.Range("B14").Value = 'Lappeenrannan kaupunki'
$ws.Range("C14").Value = 44105
$ws.Range("D14").Value = 2200
$ws.Range("E14").Value = 67179
$ws.Range("F14").Value = '2790230-3'
$ws.Range("G14").Value = 'SINGA OY'
$ws.Range("H14").Value = 45
$ws.Range("I14").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J14").Value = 4600
$ws.Range("K14").Value = 'Muu materiaali'
$ws.Range("L14").Value = 249
$ws.Range("M14").Value = 20
$ws.Range("N14").Value = 'Lapset ja nuoret'
$ws.Range("O14").Value = 222
$ws.Range("P14").Value = 'Nuorisopalvelut'
$ws.Range("Q14").Value = 2230

# Row 15
$ws.Range("A15").Value = 405
$ws.Range("B15").Value = 'Lappeenrannan kaupunki'
$ws.Range("C15").Value = 44105
$ws.Range("D15").Value = 2200
$ws.Range("E15").Value = 67179
$ws.Range("F15").Value = '2790230-3'
$ws.Range("G15").Value = 'SINGA OY'
$ws.Range("H15").Value = 45
$ws.Range("I15").Value = 'Aineet, tarvikkeet ja tavarat'
$ws.Range("J15").Value = 4600
$ws.Range("K15").Value = 'Muu materiaali'
$ws.Range("L15").Value = 12321
$ws.Range("M15").Value = 20
$ws.Range("N15").Value = 'Lapset ja nuoret'
$ws.Range("O15").Value = 222
$ws.Range("P15").Value = 'Nuorisopalvelut'
$ws.Range("Q15").Value = 2230

# Row 16
$ws.Range("A16").Value = 405
$ws.Range("B16").Value = 'Lappeenrannan kaupunki'
$ws.Range("C16").Value = 44115
$ws.Range("D16").Value = 3460
$ws.Range("E16").Value = 76785
$ws.Range("F16").Value = '2867588-6'
$ws.Range("G16").Value = 'SAIMAARIUM OY'
$ws.Range("H16").Value = 43
$ws.Range("I16").Value = 'Palvelujen ostot'
$ws.Range("J16").Value = 4440
$ws.Range("K16").Value = 'Koulutus- ja kulttuuripalvelut'
$ws.Range("L16").Value = 231
$ws.Range("M16").Value = 20
$ws.Range("N16").Value = 'Lapset ja nuoret'
$ws.Range("O16").Value = 221
$ws.Range("P16").Value = 'Lukiokoulutus'
$ws.Range("Q16").Value = 2211

# Row 17
$ws.Range("A17").Value = 405
$ws.Range("B17").Value = 'Lappeenrannan kaupunki'
$ws.Range("C17").Value = 44115
$ws.Range("D17").Value = 3460
$ws.Range("E17").Value = 76785
$ws.Range("F17").Value = '2867588-6'
$ws.Range("G17").Value = 'SAIMAARIUM OY'
$ws.Range("H17").Value = 43
$ws.Range("I17").Value = 'Palvelujen ostot'
$ws.Range("J17").Value = 4440
$ws.Range("K17").Value = 'Koulutus- ja kulttuuripalvelut'
$ws.Range("L17").Value = 50.4
$ws.Range("M17").Value = 20
$ws.Range("N17").Value = 'Lapset ja nuoret'
$ws.Range("O17").Value = 221
$ws.Range("P17").Value = 'Lukiokoulutus'
$ws.Range("Q17").Value = 2211

# Column widths for the newly-widened Tiliryhmän nro / Tiliryhmä columns
$ws.Columns.Item(8).ColumnWidth = 13.15
$ws.Columns.Item(9).ColumnWidth = 25.3

# Restore the selection to match the saved view state
$ws.Range("I21").Select()
